# Investor access workbook: insert a new "Email Enabled" column (between
# "Email" and "Cc") and mark every existing row's e-mail as enabled.
# (commit: "Made import_upload_id not nil for interests" -- the sheet
# change is the accompanying template update: new Email Enabled column.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 0. Snapshot everything that the engine will NOT auto-shift when we
#    insert a column (comments + hyperlinks stay "glued" to their old
#    cell refs instead of moving with the data), so we can rebuild it
#    in the right place afterwards.
# ---------------------------------------------------------------------

# Header comments living in columns E..J -- these need to slide one
# column to the right (E->F, F->G, G->H, H->I, I->J, J->K).
$srcCommentCols = @("E","F","G","H","I","J")
$dstCommentCols = @("F","G","H","I","J","K")
$commentTexts = @{}
foreach ($col in $srcCommentCols) {
    $cell = $ws.Range($col + "1")
    if ($cell.Comment -ne $null) {
        $commentTexts[$col] = $cell.Comment.Text()
    }
}

# Hyperlinked e-mail cells: D2:D4 stay put, E2:E4 slide to F2:F4. Grab
# the display text now (it's always the same as the mailto target) so
# we can recreate the links after the shift.
$dLinkText = @{}
$eLinkText = @{}
foreach ($r in 2..4) {
    $dLinkText[$r] = $ws.Range("D" + $r).Text
    $eLinkText[$r] = $ws.Range("E" + $r).Text
}

# ---------------------------------------------------------------------
# 1. Insert the new column. This shifts cell values/styles/column
#    widths/data validation/dimension for us.
# ---------------------------------------------------------------------
$ws.Columns("E:E").Insert()

# ---------------------------------------------------------------------
# 2. Populate the new "Email Enabled" column.
# ---------------------------------------------------------------------
$ws.Range("E1").Value = "Email Enabled"
$ws.Range("E2").Value = "Yes"
$ws.Range("E3").Value = "Yes"
$ws.Range("E4").Value = "Yes"

# ---------------------------------------------------------------------
# 3. Re-home the header comments that used to sit on E1:J1.
# ---------------------------------------------------------------------
for ($i = 0; $i -lt $srcCommentCols.Count; $i++) {
    $srcCol = $srcCommentCols[$i]
    $dstCol = $dstCommentCols[$i]
    $txt = $commentTexts[$srcCol]
    if ($txt -eq $null) { continue }
    $destCell = $ws.Range($dstCol + "1")
    if ($destCell.Comment -ne $null) {
        $destCell.Comment.Text($txt) | Out-Null
    } else {
        $destCell.AddComment($txt) | Out-Null
    }
}

# The new E1 ("Email Enabled") never carries a comment of its own -- drop
# the stale one the column-insert left glued to that ref.
if ($ws.Range("E1").Comment -ne $null) {
    $ws.Range("E1").Comment.Delete()
}

# ---------------------------------------------------------------------
# 4. Rebuild the mailto hyperlinks. Deleting a single hyperlink clears
#    the whole sheet's collection in this engine, so wipe once and
#    recreate all six in their post-shift locations.
# ---------------------------------------------------------------------
$ws.Hyperlinks.Delete()
foreach ($r in 2..4) {
    $ws.Hyperlinks.Add($ws.Range("D" + $r), ("mailto:" + $dLinkText[$r])) | Out-Null
    $ws.Range("D" + $r).Style = "Hyperlink"
    $ws.Hyperlinks.Add($ws.Range("F" + $r), ("mailto:" + $eLinkText[$r])) | Out-Null
    $ws.Range("F" + $r).Style = "Hyperlink"
}

# ---------------------------------------------------------------------
# 5. Match the saved selection state recorded in the diff.
# ---------------------------------------------------------------------
$ws.Range("E3:E4").Select() | Out-Null
